$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.137.80'
$ws.Range("E2").Value = '  +6.64%  '
$ws.Range("D3").Value = '3.565.55'
$ws.Range("E3").Value = '  +10.65%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '552.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.67%  '
$ws.Range("D7").Value = '3.558.63'
$ws.Range("E7").Value = '  +10.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.610'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.52%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.636'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.21%  '
$ws.Range("E11").Value = '  +14.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("E13").Value = '  +7.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '4.135.12'
$ws.Range("E15").Value = '  +10.76%  '
$ws.Range("D16").Value = '3.565.37'
$ws.Range("E16").Value = '  +10.82%  '
$ws.Range("E17").Value = '  +5.12%  '
$ws.Range("D18").Value = '67.082.69'
$ws.Range("E18").Value = '  +6.96%  '
$ws.Range("E19").Value = '  +6.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.996'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '434.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +18.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.19%  '
$ws.Range("E24").Value = '  +4.13%  '
$ws.Range("E25").Value = '  +4.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("E27").Value = '  +9.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '649.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("E33").Value = '  +3.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.53%  '
$ws.Range("E35").Value = '  +6.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.79'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.94%  '
$ws.Range("E37").Value = '  +24.61%  '
$ws.Range("D38").Value = '0.0₃0827'
$ws.Range("E38").Value = '  +17.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.392'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.14%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '3.035.84'
$ws.Range("E45").Value = '  +5.85%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.88'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.78%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.14%  '
$ws.Range("E48").Value = '  +7.33%  '
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("E50").Value = '  +5.14%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.99%  '
